$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 7 de Octubre de 2020 a las 11:05"

# Update country name labels (column A) for rows whose country order changed
$ws.Cells.Item(22, 1).Value = "Filipinas"
$ws.Cells.Item(23, 1).Value = "Turquia"

$ws.Cells.Item(82, 1).Value = "Tunez"
$ws.Cells.Item(83, 1).Value = "Corea del Sur"
$ws.Cells.Item(84, 1).Value = "Bulgaria"

$ws.Cells.Item(136, 1).Value = "Sri Lanka"
$ws.Cells.Item(137, 1).Value = "Reunion"

$ws.Cells.Item(154, 1).Value = "Letonia"
$ws.Cells.Item(155, 1).Value = "Belice"
$ws.Cells.Item(156, 1).Value = "Polinesia Francesa"
$ws.Cells.Item(157, 1).Value = "Burkina Faso"

$ws.Cells.Item(207, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(208, 1).Value = "Santa Lucia"

# Update numeric data (columns B-H) for changed rows
# Row 22
$ws.Cells.Item(22, 2).Value = 329637
$ws.Cells.Item(22, 3).Value = 2825
$ws.Cells.Item(22, 4).Value = 273723
$ws.Cells.Item(22, 5).Value = 49989
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(22, 7).Value = 60
$ws.Cells.Item(22, 8).Value = 5925

# Row 23
$ws.Cells.Item(23, 2).Value = 327557
$ws.Cells.Item(23, 3).Value = 0
$ws.Cells.Item(23, 4).Value = 287599
$ws.Cells.Item(23, 5).Value = 31405
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(23, 8).Value = 8553

# Row 25
$ws.Cells.Item(25, 2).Value = 315714
$ws.Cells.Item(25, 3).Value = 4538
$ws.Cells.Item(25, 4).Value = 240291
$ws.Cells.Item(25, 5).Value = 63951
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 98
$ws.Cells.Item(25, 8).Value = 11472

# Row 41
$ws.Cells.Item(41, 2).Value = 107319
$ws.Cells.Item(41, 3).Value = 3003
$ws.Cells.Item(41, 4).Value = 75346
$ws.Cells.Item(41, 5).Value = 29181
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(41, 7).Value = 75
$ws.Cells.Item(41, 8).Value = 2792

# Row 65
$ws.Cells.Item(65, 2).Value = 50848
$ws.Cells.Item(65, 3).Value = 1029
$ws.Cells.Item(65, 4).Value = 40499
$ws.Cells.Item(65, 5).Value = 9519
$ws.Cells.Item(65, 6).Value = 0
$ws.Cells.Item(65, 7).Value = 8
$ws.Cells.Item(65, 8).Value = 830

# Row 82
$ws.Cells.Item(82, 2).Value = 24542
$ws.Cells.Item(82, 3).Value = 2312
$ws.Cells.Item(82, 4).Value = 5032
$ws.Cells.Item(82, 5).Value = 19146
$ws.Cells.Item(82, 6).Value = 0
$ws.Cells.Item(82, 7).Value = 43
$ws.Cells.Item(82, 8).Value = 364

# Row 83
$ws.Cells.Item(83, 2).Value = 24353
$ws.Cells.Item(83, 3).Value = 114
$ws.Cells.Item(83, 4).Value = 22334
$ws.Cells.Item(83, 5).Value = 1594
$ws.Cells.Item(83, 6).Value = 0
$ws.Cells.Item(83, 7).Value = 3
$ws.Cells.Item(83, 8).Value = 425

# Row 84
$ws.Cells.Item(84, 2).Value = 22306
$ws.Cells.Item(84, 3).Value = 0
$ws.Cells.Item(84, 4).Value = 15310
$ws.Cells.Item(84, 5).Value = 6134
$ws.Cells.Item(84, 6).Value = 0
$ws.Cells.Item(84, 8).Value = 862

# Row 91
$ws.Cells.Item(91, 2).Value = 18447
$ws.Cells.Item(91, 3).Value = 363
$ws.Cells.Item(91, 4).Value = 16308
$ws.Cells.Item(91, 5).Value = 1830
$ws.Cells.Item(91, 6).Value = 0
$ws.Cells.Item(91, 7).Value = 5
$ws.Cells.Item(91, 8).Value = 309

# Row 136
$ws.Cells.Item(136, 2).Value = 4442
$ws.Cells.Item(136, 3).Value = 190
$ws.Cells.Item(136, 4).Value = 3274
$ws.Cells.Item(136, 5).Value = 1155
$ws.Cells.Item(136, 6).Value = 0
$ws.Cells.Item(136, 8).Value = 13

# Row 137
$ws.Cells.Item(137, 2).Value = 4328
$ws.Cells.Item(137, 3).Value = 0
$ws.Cells.Item(137, 4).Value = 3360
$ws.Cells.Item(137, 5).Value = 952
$ws.Cells.Item(137, 6).Value = 0
$ws.Cells.Item(137, 8).Value = 16

# Row 141
$ws.Cells.Item(141, 2).Value = 3684
$ws.Cells.Item(141, 3).Value = 25
$ws.Cells.Item(141, 4).Value = 2813
$ws.Cells.Item(141, 5).Value = 804
$ws.Cells.Item(141, 6).Value = 0

# Row 154
$ws.Cells.Item(154, 2).Value = 2261
$ws.Cells.Item(154, 3).Value = 67
$ws.Cells.Item(154, 4).Value = 1322
$ws.Cells.Item(154, 5).Value = 899
$ws.Cells.Item(154, 6).Value = 0
$ws.Cells.Item(154, 7).Value = 0
$ws.Cells.Item(154, 8).Value = 40

# Row 155
$ws.Cells.Item(155, 2).Value = 2243
$ws.Cells.Item(155, 3).Value = 39
$ws.Cells.Item(155, 4).Value = 1392
$ws.Cells.Item(155, 5).Value = 817
$ws.Cells.Item(155, 6).Value = 0
$ws.Cells.Item(155, 7).Value = 4
$ws.Cells.Item(155, 8).Value = 34

# Row 156
$ws.Cells.Item(156, 2).Value = 2228
$ws.Cells.Item(156, 3).Value = 0
$ws.Cells.Item(156, 4).Value = 1769
$ws.Cells.Item(156, 5).Value = 450
$ws.Cells.Item(156, 6).Value = 0
$ws.Cells.Item(156, 8).Value = 9

# Row 157
$ws.Cells.Item(157, 2).Value = 2197
$ws.Cells.Item(157, 3).Value = 0
$ws.Cells.Item(157, 4).Value = 1441
$ws.Cells.Item(157, 5).Value = 697
$ws.Cells.Item(157, 6).Value = 0
$ws.Cells.Item(157, 8).Value = 59

